$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Two new leading "aggregate" columns are inserted (All / Europe), the
# existing country columns are re-ordered (Europe first, then the rest)
# and "United States" is renamed to "USA".
# (A1 is left untouched - it is already blank.)
$ws.Range("B1").Value = "`$ bold('All')"
$ws.Range("C1").Value = "`$ bold('Europe')"
$ws.Range("D1").Value = "France"
$ws.Range("E1").Value = "Germany"
$ws.Range("F1").Value = "Italy"
$ws.Range("G1").Value = "Poland"
$ws.Range("H1").Value = "Spain"
$ws.Range("I1").Value = "United Kingdom"
$ws.Range("J1").Value = "Switzerland"
$ws.Range("K1").Value = "Japan"
$ws.Range("L1").Value = "Saudi Arabia"
$ws.Range("M1").Value = "USA"

# --- Row 2 data ---------------------------------------------------------
$ws.Range("B2").Value = 0.548749204972874
$ws.Range("C2").Value = 0.612606004275778
$ws.Range("D2").Value = 0.616724928283738
$ws.Range("E2").Value = 0.620808072627315
$ws.Range("F2").Value = 0.75031578480093
$ws.Range("G2").Value = 0.495076291993524
$ws.Range("H2").Value = 0.6117343914093
$ws.Range("I2").Value = 0.554433842033875
$ws.Range("J2").Value = 0.529872342725065
$ws.Range("K2").Value = 0.438153093874799
$ws.Range("L2").Value = 0.677827112481047
$ws.Range("M2").Value = 0.508222265196169

# --- Row 3 data ---------------------------------------------------------
$ws.Range("B3").Value = 0.488009346515533
$ws.Range("C3").Value = 0.55532380171867
$ws.Range("D3").Value = 0.588306427984715
$ws.Range("E3").Value = 0.528389866344171
$ws.Range("F3").Value = 0.602852192238351
$ws.Range("G3").Value = 0.547884560154944
$ws.Range("H3").Value = 0.567216620094722
$ws.Range("I3").Value = 0.542572051573976
$ws.Range("J3").Value = 0.360111742646164
$ws.Range("K3").Value = 0.351009115929929
$ws.Range("L3").Value = 0.671355992475847
$ws.Range("M3").Value = 0.447291263751444
